$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 531
$ws.Range("F5").Value = 2541
$ws.Range("F7").Value = 92
$ws.Range("F9").Value = 1718
$ws.Range("F10").Value = 1718
$ws.Range("F11").Value = 1429
$ws.Range("F13").Value = 1463
$ws.Range("F14").Value = 26
$ws.Range("F16").Value = 1009
$ws.Range("F19").Value = 259
$ws.Range("F20").Value = 7632
$ws.Range("F21").Value = 8730
$ws.Range("F22").Value = 61
$ws.Range("F24").Value = 432
$ws.Range("F26").Value = 104
$ws.Range("F27").Value = 276
$ws.Range("F32").Value = 1539
$ws.Range("F34").Value = 272
$ws.Range("F39").Value = 826
$ws.Range("F42").Value = 379
$ws.Range("F46").Value = 230
$ws.Range("F47").Value = 8

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 25
$ws.Range("F3").Value = 12
$ws.Range("F4").Value = 79
$ws.Range("F19").Value = 329

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 314

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 25
$ws.Range("F6").Value = 314
$ws.Range("F8").Value = 531
$ws.Range("F9").Value = 2541
$ws.Range("F11").Value = 92
$ws.Range("F13").Value = 1718
$ws.Range("F14").Value = 1718
$ws.Range("F16").Value = 1463
$ws.Range("F18").Value = 1009
$ws.Range("F20").Value = 12
$ws.Range("F22").Value = 79
$ws.Range("F23").Value = 259
$ws.Range("F24").Value = 7632
$ws.Range("F25").Value = 7632
$ws.Range("F26").Value = 8730
$ws.Range("F27").Value = 61
$ws.Range("F29").Value = 432
$ws.Range("F30").Value = 104
$ws.Range("F31").Value = 276
$ws.Range("F34").Value = 1539
$ws.Range("F36").Value = 272
$ws.Range("F41").Value = 826
$ws.Range("F43").Value = 379
$ws.Range("F47").Value = 230
$ws.Range("F48").Value = 8
$ws.Range("F50").Value = 329
